# Apply the "Add files via upload" update to the IK Konj+Destatis_HWWI data sheet.
# The underlying change is: the 2025-Q2 row (row 43) gets its quarterly figures
# filled in (columns C:F and O:W), and the still-empty 2025-Q4 placeholder row
# (row 45) is removed, shrinking the used range from A1:W45 to A1:W44.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the newly-reported figures for 2025 Q2 (row 43) ---
$ws.Range("C43").Value = 206.45
$ws.Range("D43").Value = 152.44
$ws.Range("E43").Value = 155.37
$ws.Range("F43").Value = 376.9

$ws.Range("O43").Value = 509
$ws.Range("P43").Value = 90306.666666666686
$ws.Range("Q43").Value = 34012667
$ws.Range("R43").Value = 1178328666
$ws.Range("S43").Value = 6355352000
$ws.Range("T43").Value = 3367534333
$ws.Range("U43").Value = 2987817666
$ws.Range("V43").Value = 1757211000
$ws.Range("W43").Value = 1230606667

# --- Drop the still-empty 2025 Q4 row entirely (used range shrinks to row 44) ---
$ws.Rows("45").Delete()
